# The workbook's active sheet (Feuille2) is not the sheet that needs edits -
# the target data lives on the worksheet literally named "sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 4 new columns starting at column EW (153), pushing the existing
# "Adresse de courriel" / empty columns (formerly EW/EX) to FA/FB. Excel's
# column insert naturally copies formatting from the column to the left,
# which reproduces style id 3 on the new cells and keeps style id 5 on the
# shifted email/empty columns.
$ws.Range("EW1:EZ9").EntireColumn.Insert()

# Populate the newly inserted EW:EZ columns by repeating the same cyclical
# "Alain / Henri / Tony / Dulcinee" (row 1) and "OUI/NON" (rows 2-9) pattern
# that already fills columns E:EV, i.e. copy each row's E:H values into the
# new EW:EZ cells.
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 153).Value2 = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 154).Value2 = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 155).Value2 = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 156).Value2 = $ws.Cells.Item($r, 8).Value2
}
